$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 823
$ws.Range("I15").Value = 823
$ws.Range("K15").Value = 2469
$ws.Range("M15").Value = -2300

# Row 98
$ws.Range("H98").Value = 602.9091
$ws.Range("I98").Value = 516.8889
$ws.Range("J98").Value = 990
$ws.Range("K98").Value = 516.8889
$ws.Range("L98").Value = 990
$ws.Range("M98").Value = 981.1111
$ws.Range("N98").Value = -3986

# Row 122
$ws.Range("H122").Value = 602.9091
$ws.Range("I122").Value = 516.8889
$ws.Range("J122").Value = 990
$ws.Range("K122").Value = 1550.6667
$ws.Range("L122").Value = 2970
$ws.Range("M122").Value = 899.3332999999998
$ws.Range("N122").Value = -7870

# Row 135
$ws.Range("H135").Value = 2320.0952
$ws.Range("I135").Value = 2264.3157
$ws.Range("K135").Value = 20378.8413
$ws.Range("M135").Value = -17843.8413

# Row 137
$ws.Range("H137").Value = 5299
$ws.Range("I137").Value = 5299
$ws.Range("K137").Value = 15897
$ws.Range("M137").Value = -13347

# Row 138
$ws.Range("H138").Value = 3000
$ws.Range("J138").Value = 3000
$ws.Range("L138").Value = 9000
$ws.Range("N138").Value = -19280

# Row 141
$ws.Range("H141").Value = 12646
$ws.Range("I141").Value = 16728.334
$ws.Range("J141").Value = 399
$ws.Range("K141").Value = 50185.00199999999
$ws.Range("L141").Value = 1197
$ws.Range("M141").Value = -45005.00199999999
$ws.Range("N141").Value = -11557

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

# Row 40
$ws.Range("H40").Value = 200000
$ws.Range("J40").Value = 200000
$ws.Range("L40").Value = 200000
$ws.Range("N40").Value = -200530

# Row 103
$ws.Range("H103").Value = 22147.5
$ws.Range("J103").Value = 22147.5
$ws.Range("L103").Value = 22147.5
$ws.Range("N103").Value = -24491.5

# Row 107
$ws.Range("H107").Value = 3000
$ws.Range("I107").Value = 3000
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 3000
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -1080

# Row 134
$ws.Range("H134").Value = 1517.375
$ws.Range("I134").Value = 1034.1428
$ws.Range("J134").Value = 4900
$ws.Range("K134").Value = 3102.4284
$ws.Range("L134").Value = 14700
$ws.Range("M134").Value = -567.4284000000002
$ws.Range("N134").Value = -19770

$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 1000
$ws.Range("I4").Value = 1000
$ws.Range("K4").Value = 1000
$ws.Range("M4").Value = -888

# Row 5
$ws.Range("H5").Value = 374.85715
$ws.Range("I5").Value = 354
$ws.Range("J5").Value = 500
$ws.Range("K5").Value = 354
$ws.Range("L5").Value = 500
$ws.Range("M5").Value = -242
$ws.Range("N5").Value = -724

# Row 12
$ws.Range("H12").Value = 1346.75
$ws.Range("J12").Value = 2000
$ws.Range("L12").Value = 2000
$ws.Range("N12").Value = -2340

# Row 25
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()

# Row 35
$ws.Range("H35").Value = 2333.8572
$ws.Range("J35").Value = 4250
$ws.Range("L35").Value = 4250
$ws.Range("N35").Value = -4838

# Row 132
$ws.Range("H132").Value = 8048
$ws.Range("I132").Value = 8048
$ws.Range("K132").Value = 24144
$ws.Range("M132").Value = -21614

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 500
$ws.Range("I5").Value = 500
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1500
$ws.Range("L5").ClearContents()
$ws.Range("M5").Value = -1388
$ws.Range("N5").Value = 0

# Row 37
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").ClearContents()
$ws.Range("N37").Value = 0

# Row 107
$ws.Range("H107").Value = 1005.8889
$ws.Range("I107").Value = 1287.6
$ws.Range("J107").Value = 653.75
$ws.Range("K107").Value = 3862.8
$ws.Range("L107").Value = 1961.25
$ws.Range("M107").Value = -1942.8
$ws.Range("N107").Value = -5801.25

# Row 117
$ws.Range("H117").Value = 6377.75
$ws.Range("I117").Value = 614.5
$ws.Range("K117").Value = 1843.5
$ws.Range("M117").Value = 1598.5

# Row 129
$ws.Range("H129").Value = 12856.5
$ws.Range("J129").Value = 24883
$ws.Range("L129").Value = 74649
$ws.Range("N129").Value = -84649

# Row 131
$ws.Range("H131").Value = 2469.0833
$ws.Range("I131").Value = 1604.8334
$ws.Range("J131").Value = 3333.3333
$ws.Range("K131").Value = 4814.5002
$ws.Range("L131").Value = 9999.999899999999
$ws.Range("M131").Value = 225.4997999999996
$ws.Range("N131").Value = -20079.9999

# Row 135
$ws.Range("H135").Value = 500
$ws.Range("I135").Value = 500
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 4500
$ws.Range("L135").ClearContents()
$ws.Range("M135").Value = -1965
$ws.Range("N135").Value = 0

$ws = $wb.Worksheets.Item("GSM")
# Row 24
$ws.Range("H24").Value = 17801.4
$ws.Range("J24").Value = 17801.4
$ws.Range("L24").Value = 17801.4
$ws.Range("N24").Value = -18147.4

# Row 70
$ws.Range("H70").Value = 19498.334
$ws.Range("I70").Value = 3500
$ws.Range("K70").Value = 3500
$ws.Range("M70").Value = -3230

# Row 73
$ws.Range("H73").Value = 19498.334
$ws.Range("I73").Value = 3500
$ws.Range("K73").Value = 3500
$ws.Range("M73").Value = -2564

# Row 132
$ws.Range("H132").Value = 1585.7142
$ws.Range("I132").Value = 1516.6666
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 4549.9998
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -2019.9998
$ws.Range("N132").Value = -11060

$ws = $wb.Worksheets.Item("LTW")
# Row 42
$ws.Range("H42").Value = 1000
$ws.Range("I42").Value = 1000
$ws.Range("J42").Value = 1000
$ws.Range("K42").Value = 1000
$ws.Range("L42").Value = 1000
$ws.Range("M42").Value = -437
$ws.Range("N42").Value = -2126

# Row 49
$ws.Range("H49").Value = 1000
$ws.Range("I49").Value = 1000
$ws.Range("J49").Value = 1000
$ws.Range("K49").Value = 1000
$ws.Range("L49").Value = 1000
$ws.Range("M49").Value = -853
$ws.Range("N49").Value = -1294

# Row 55
$ws.Range("H55").Value = 413.57144
$ws.Range("I55").Value = 413.57144
$ws.Range("K55").Value = 413.57144
$ws.Range("M55").Value = -240.57144

# Row 93
$ws.Range("H93").Value = 849.6667
$ws.Range("I93").Value = 849.6667
$ws.Range("K93").Value = 849.6667
$ws.Range("M93").Value = 398.3333

# Row 136
$ws.Range("H136").Value = 67917
$ws.Range("I136").Value = 26973
$ws.Range("K136").Value = 80919
$ws.Range("M136").Value = -78369

$ws = $wb.Worksheets.Item("WVR")
# Row 21
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

# Row 35
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()

# Row 122
$ws.Range("H122").Value = 850
$ws.Range("J122").Value = 1000
$ws.Range("L122").Value = 3000
$ws.Range("N122").Value = -7900

# Row 132
$ws.Range("H132").Value = 3998
$ws.Range("I132").Value = 3997
$ws.Range("K132").Value = 11991
$ws.Range("M132").Value = -9461
